$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 73, shifting existing rows 73:107 down to 74:108
$ws.Rows.Item(73).Insert()

# Copy formatting/style for the date cell from the row below (already shifted) so the
# new row keeps the same date number format used throughout column D
$ws.Range("D74").Copy()
$ws.Range("D73").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Populate the new row 73 with its data (columns that stay the same as the rest of the
# dataset for this market are filled too, matching the pattern of every other row)
$ws.Range("A73").Value = 7
$ws.Range("B73").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C73").Value = "Ñuble"
$ws.Range("D73").Value = 45141
$ws.Range("E73").Value = 16
$ws.Range("F73").Value = 100112001
$ws.Range("G73").Value = "Berenjena"
$ws.Range("H73").Value = "Sin especificar"
$ws.Range("I73").Value = "Primera"
$ws.Range("J73").Value = 30
$ws.Range("K73").Value = 9000
$ws.Range("L73").Value = 9000
$ws.Range("M73").Value = 9000
$ws.Range("N73").Value = "$/caja 60 unidades"
$ws.Range("O73").Value = "Región de Arica y Parinacota"
$ws.Range("P73").Value = 150
$ws.Range("Q73").Value = 60
$ws.Range("R73").Value = "Hortaliza"
